$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title: "Fantasy Points for NFL Players" -> "...Quarterbacks"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Fantasy Points for NFL Players", $true, $false, $false, $false, $false, $true, 1, $false, "Fantasy Points for NFL Quarterbacks", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Typo fix: "Thos who enjoy" -> "Those who enjoy" (insert the
#    missing "e" right after "Thos").
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Thos who enjoy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$thosEnd = $r.Start + 4                      # position right after "Thos"
$insPoint = $d.Range($thosEnd, $thosEnd)
$insPoint.InsertAfter("e")
$goBackPos = $thosEnd + 1                    # right after the inserted "e" ("Those|")

# ---------------------------------------------------------------------
# 3) Relocate the hidden "_GoBack" bookmark from its old spot (right
#    after "ant. ", before "Since we are focusing...") to its new spot
#    (right after "Those", before " who enjoy...").
# ---------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$newBookmarkRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange) | Out-Null
